$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 402
$ws.Range("C2").Value = 554

$ws.Range("B3").Value = 465
$ws.Range("C3").Value = 479

$ws.Range("B4").Value = 257
$ws.Range("C4").Value = 305

$ws.Range("B6").Value = 183
$ws.Range("C6").Value = 250

$ws.Range("A7").Value = "Serie A"
$ws.Range("B7").Value = 291
$ws.Range("C7").Value = 444

$ws.Range("B7").Select()
